$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the PREÇO and DATA ATUAL columns (B and C), including header and data
$ws.Range("B1:C14").Clear()

# Add the new product rows in column A
$ws.Range("A15").Value = "Aluminio"
$ws.Range("A16").Value = "Niquel"
$ws.Range("A17").Value = "Zinco"
$ws.Range("A18").Value = "Titanio"
$ws.Range("A19").Value = "Cacau"
